$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G3").Value = 60
$ws1.Range("F5").Value = 7634
$ws1.Range("F6").Value = 5553
$ws1.Range("F7").Value = 455
$ws1.Range("F12").Value = 186

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 83

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G3").Value = 60
$ws4.Range("F5").Value = 7634
$ws4.Range("F6").Value = 5553
$ws4.Range("F7").Value = 455
$ws4.Range("F12").Value = 83
$ws4.Range("F14").Value = 186
